$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'Rodovia Raposo Tavares'
$ws.Range('A3').Value = 'Rua Hermenegildo D''Andréa'
$ws.Range('A4').Value = 'Avenida das Letras'
$ws.Range('A5').Value = 'Rua Clara'
$ws.Range('A6').Value = 'Avenida Moysés Sayão'
$ws.Range('A7').Value = 'Alameda Francisco Cristófani'
$ws.Range('A8').Value = ''
$ws.Range('A9').Value = 'Rua Adelmar Tavares'
$ws.Range('A10').Value = 'Rua Icó'
$ws.Range('A11').Value = 'Travessa Herval D''Oeste'
$ws.Range('A12').Value = 'Rua Pernambucana'
$ws.Range('A13').Value = 'Rua Ferroviária Bráulio dos Reis'
$ws.Range('A14').Value = 'Rua Adolfo Rodrigues'
$ws.Range('A15').Value = 'Avenida Virgilio Cardoso Pinna'
$ws.Range('A16').Value = 'Rua Ancara'
$ws.Range('A17').Value = 'Rua Gamboa'
$ws.Range('A18').Value = 'Rua dos Cravos'
$ws.Range('A19').Value = 'Rua Oslo'
$ws.Range('A20').Value = 'Rua Anhanga'
$ws.Range('A21').Value = 'Rua Vicentina'
$ws.Range('A22').Value = 'Rua Bernardo da Veiga'
$ws.Range('A23').Value = 'Avenida Doutor Ricardo Jafet'
$ws.Range('A24').Value = ''
$ws.Range('A25').Value = 'Rua Júlio Verne'
$ws.Range('A26').Value = 'Rua Vasco Fernandes Coutinho'
$ws.Range('A27').Value = 'Rua Coroados'
$ws.Range('A28').Value = 'Rua José Yazigi'
$ws.Range('A29').Value = 'Rua Refinaria Mataripe'
$ws.Range('A30').Value = 'Rua Odete'
$ws.Range('A31').Value = ''
$ws.Range('A32').Value = 'Avenida Nova Independência'
$ws.Range('A33').Value = 'Rua Bartolomeu Caporali'
$ws.Range('A34').Value = 'Avenida Dom Pedro I'
$ws.Range('A35').Value = 'Avenida Valentim Magalhães'
$ws.Range('A36').Value = 'Rua das Garças'
$ws.Range('A37').Value = 'Avenida Araucária'
$ws.Range('A38').Value = 'Rua Nilo Peçanha'
$ws.Range('A39').Value = 'Avenida Vereador Carlito Cordeiro'
$ws.Range('A40').Value = 'Rua Itanhomi'
$ws.Range('A41').Value = ''
$ws.Range('A42').Value = 'Rua Rogério Giorgi'
$ws.Range('A43').Value = 'Rua Rio Grande do Sul'
$ws.Range('A44').Value = 'Rua Bernardo da Veiga'
$ws.Range('A45').Value = 'Rua Werner Goldberg'
$ws.Range('A46').Value = 'Rua Clélia'
$ws.Range('A47').Value = 'Rua Sebastião Annunciatto'
$ws.Range('A48').Value = 'Avenida Professor Luiz Ignácio Anhaia Mello'
$ws.Range('A49').Value = 'Estrada Geral Vargem Pequena'
$ws.Range('A50').Value = ''
$ws.Range('A51').Value = 'Avenida Firestone'
$ws.Range('A52').Value = 'Rua Caquito'
$ws.Range('A53').Value = 'Rua Doutor Virgílio de Carvalho Pinto'
$ws.Range('A54').Value = 'Rua Cavour'
$ws.Range('A55').Value = ''
$ws.Range('A56').Value = 'Rua Tangânica'
$ws.Range('A57').Value = ''
$ws.Range('A58').Value = ''
$ws.Range('A59').Value = 'Rua Ministro José Geraldo Rodrigues Alkmin'
$ws.Range('A60').Value = 'Travessa Apeninos'
$ws.Range('A61').Value = 'Rua Refinaria Mataripe'
$ws.Range('A62').Value = 'Rua São João'
$ws.Range('A63').Value = 'Rua 30 de Outubro'
$ws.Range('A64').Value = 'Rua Ibitirama'
$ws.Range('A65').Value = 'Rua Coréia'
$ws.Range('A66').Value = 'Rua Alberto de Oliveira'
$ws.Range('A67').Value = 'Rua Raul Pompéia'
$ws.Range('A68').Value = ''
$ws.Range('A69').Value = 'Rua Adriano Theodósio Serra'
$ws.Range('A70').Value = 'Rua Conde Juliano'
$ws.Range('A71').Value = 'Alameda Princeza Januária'
$ws.Range('A72').Value = 'Rua Abolição'
$ws.Range('A73').Value = 'Rua Domingos Garcia Velho'
$ws.Range('A74').Value = 'Rua Piracanjuba'
$ws.Range('A75').Value = 'Rua Quintino Bocaiúva'
$ws.Range('A76').Value = ''
$ws.Range('A77').Value = 'Rua Daniel Berg'
$ws.Range('A78').Value = 'Rodovia BR-262'
$ws.Range('A79').Value = 'Estrada Manoel Lages do Chao'
$ws.Range('A80').Value = 'Rua Onze de Agosto'
$ws.Range('A81').Value = 'Rua Campo Grande'
$ws.Range('A82').Value = 'Rua Augusto Maass'
$ws.Range('A83').Value = ''
$ws.Range('A84').Value = ''
$ws.Range('A85').Value = 'Rua Leopoldo Schmidt'
$ws.Range('A86').Value = ''
$ws.Range('A87').Value = 'Rua Jandiro Joaquim Pereira'
$ws.Range('A88').Value = 'Rua João Maluf'
$ws.Range('A89').Value = 'Travessa João Mendes'
$ws.Range('A90').Value = 'Rua Alfredo Faria de Souza'
$ws.Range('A91').Value = 'Avenida Mendes da Rocha'
$ws.Range('A92').Value = ''
$ws.Range('A93').Value = ''
$ws.Range('A94').Value = 'Rua Quintino Bocaiúva'
$ws.Range('A95').Value = 'Rua Eugênia de Carvalho'
$ws.Range('A96').Value = 'Rua Álvaro Lins'
$ws.Range('A97').Value = 'Rua Doutor Fleury Silveira'
$ws.Range('A98').Value = 'Rua PP 8'
$ws.Range('A99').Value = 'Rua Kalil Filho'
$ws.Range('A100').Value = 'Rua CJ 04'
$ws.Range('A101').Value = 'Rua José de Oliveira Coelho'
$ws.Range('A102').Value = 'Rua das Azaléas'
$ws.Range('A103').Value = 'Avenida Aruanã'
$ws.Range('A104').Value = 'Travessa Cantareira'
$ws.Range('A105').Value = 'Rua Antônio Salviano de Rezende'
$ws.Range('A106').Value = 'Avenida Francisco José Resende'
